# Updated Features Roadmap for IDA
# Adds 5 new feature-roadmap rows (S.No. 151-155) to the MOSIP_Feature_Roadmap sheet,
# covering JIRA tickets MOS-21582, MOS-21583, MOS-21584, MOS-21585 and MOS-21327.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MOSIP_Feature_Roadmap")

# ---------------------------------------------------------------------------
# Helper: stamp a single data row (A..O) using the same visual style already
# used by the existing data rows immediately above (row 153) so the new rows
# blend in with the rest of the table. NOTE: this embedded PowerShell dialect
# does not bind named (-Param) arguments inside function bodies, so we use
# positional parameters only.
# ---------------------------------------------------------------------------
function Set-RoadmapRow($Row, $SNo, $Jira, $IdentifiedDate, $Reference, $Module, $NewOrDescoped, $ChangeDescription, $Version, $ApprovalStatus, $Approver, $ApprovedDate) {

    $rowRange = "A" + $Row + ":U" + $Row

    # Bring over the formatting (borders/fonts/alignment/number formats) of the
    # last fully populated row (153) before we overwrite the cell contents.
    $ws.Range("A153:U153").Copy()
    $ws.Range($rowRange).PasteSpecial(-4122)

    $ws.Range("A" + $Row).Value = $SNo
    $ws.Range("B" + $Row).Value = $Jira
    $ws.Range("C" + $Row).Value = $IdentifiedDate
    $ws.Range("C" + $Row).NumberFormat = "d-mmm-yy"
    $ws.Range("D" + $Row).Value = $Reference
    $ws.Range("E" + $Row).Value = $Module
    $ws.Range("F" + $Row).Value = $NewOrDescoped
    $ws.Range("G" + $Row).Value = $ChangeDescription
    $ws.Range("L" + $Row).Value = $Version
    $ws.Range("M" + $Row).Value = $ApprovalStatus
    $ws.Range("N" + $Row).Value = $Approver

    # Column O (approval date) uses the same style as the neighbouring
    # Approver/Version cells, but with a date number format applied.
    $oCell = "O" + $Row
    $ws.Range("S85").Copy()
    $ws.Range($oCell).PasteSpecial(-4122)
    $ws.Range($oCell).Value = $ApprovedDate
    $ws.Range($oCell).NumberFormat = "d-mmm-yy"
}

Set-RoadmapRow 154 151 "MOS-21582" 43550 "API Specification Changes for IDA based on MDS review by Sasi/Ramesh" "ID-Authentication" "New" "Additional or Modification of attributes in API Specs based on review " 1 "Approved" "Ramesh" 43550

Set-RoadmapRow 155 152 "MOS-21583" 43550 "Design Change of ID-Repo based on Security review by Sasi/Ramesh" "ID-Authentication" "New" "Design Change of ID-Repo based on Security review by Sasi/Ramesh" 1 "Approved" "Ramesh" 43550

Set-RoadmapRow 156 153 "MOS-21584" 43550 "Design Change of IDA based on Security review by Sasi/Ramesh" "ID-Authentication" "New" "Design Change of IDA based on Security review by Sasi/Ramesh" 1 "Approved" "Ramesh" 43550

Set-RoadmapRow 157 154 "MOS-21585" 43556 "Mapping of platform address attributes in IDA based on Morrocco Address Structure" "ID-Authentication" "New" "Mapping of platform address attributes in IDA based on Morrocco Address Structure" 1 "Approved" "Shrikant" 43556

Set-RoadmapRow 158 155 "MOS-21327" 43552 "Integrate with new VID Generator API" "ID-Authentication" "New" "Integrate with the new VID generator component based on the VID policy/type defined" 1 "Approved" "Ramesh" 43552

$ws.Range("A152").Select()
